# FESTIVOS.xlsx maintenance edit
#
# 1) Rename the "BA" (Baleares) sheet to "IB" (Islas Baleares) — the sheet
#    keeps its position (5th tab) and sheetId/relationship, only the
#    display name changes.
# 2) Switch the active/selected tab from "ME" (last sheet, index 25) to
#    "CT" (index 15), which also moves workbookView.activeTab 25 -> 15 and
#    moves sheetView tabSelected="1" from ME's sheet to CT's sheet.

$wb = $excel.ActiveWorkbook

# --- 1) Rename BA -> IB ------------------------------------------------
$wb.Worksheets.Item("BA").Name = "IB"

# --- 2) Make CT the active/selected sheet ------------------------------
$wb.Worksheets.Item("CT").Activate()
$wb.Worksheets.Item("CT").Select()
